$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.3073112269268004, 1.4662627121157903, -0.7533565482959987, 2.366143232213453)"
$ws.Range("C2").Value = "JSU(-1.0529497042392317, 1.343859859563398, 3.4928454841736345, 6.545793970744947)"
$ws.Range("D2").Value = "NIG(0.9916226789776033, 0.7765575776214685, 0.7551801201722991, 3.161456260717257)"
$ws.Range("E2").Value = "EXN(3.094891655708858, 2.1031943695210895, 3.146090747745311)"
